# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" quarter sheet (inserted right after "总计") and
# prepends a matching summary row at the top of the "总计" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item(1)        # "总计"

# --- 1. Create the new "2022-Q4" worksheet right before the current first
#        quarter sheet, so it becomes the new second tab -------------------
$newSheet = $wb.Worksheets.Add($null, $summarySheet)
$newSheet.Name = "2022-Q4"

# Look the old first-quarter sheet back up by name now that the tab order
# has shifted (index-based refs would now point at the new sheet).
$firstQuarterSheet = $wb.Worksheets.Item("2022-Q3")

# Clone the header-row formatting (bold + border, centered) from the
# existing quarter sheet, and the data-row formatting (index column style)
# tiled down across all 13 data rows we are about to populate.
$firstQuarterSheet.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)

$firstQuarterSheet.Range("A2:H2").Copy()
$newSheet.Range("A2:H14").PasteSpecial(-4122)

# Fund code / scale / weight / value columns are stored as text (to keep
# leading zeros such as "001678" and fixed-decimal strings like "4.20").
$newSheet.Range("B2:G14").NumberFormat = "@"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "001678"
$newSheet.Range("C2").Value = "英大国企改革主题股票"
$newSheet.Range("D2").Value = "4.20"
$newSheet.Range("E2").Value = "92.20"
$newSheet.Range("F2").Value = "8.67"
$newSheet.Range("G2").Value = "0.3641"
$newSheet.Range("H2").Value = 1
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "003713"
$newSheet.Range("C3").Value = "英大睿盛灵活配置混合A"
$newSheet.Range("D3").Value = "2.39"
$newSheet.Range("E3").Value = "93.29"
$newSheet.Range("F3").Value = "9.28"
$newSheet.Range("G3").Value = "0.2218"
$newSheet.Range("H3").Value = 2
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "003714"
$newSheet.Range("C4").Value = "英大睿盛灵活配置混合C"
$newSheet.Range("D4").Value = "2.39"
$newSheet.Range("E4").Value = "93.29"
$newSheet.Range("F4").Value = "9.28"
$newSheet.Range("G4").Value = "0.2218"
$newSheet.Range("H4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "012202"
$newSheet.Range("C5").Value = "中加消费优选混合A"
$newSheet.Range("D5").Value = "3.65"
$newSheet.Range("E5").Value = "88.36"
$newSheet.Range("F5").Value = "3.91"
$newSheet.Range("G5").Value = "0.1427"
$newSheet.Range("H5").Value = 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "001607"
$newSheet.Range("C6").Value = "英大策略优选混合A"
$newSheet.Range("D6").Value = "0.59"
$newSheet.Range("E6").Value = "93.12"
$newSheet.Range("F6").Value = "6.80"
$newSheet.Range("G6").Value = "0.0401"
$newSheet.Range("H6").Value = 4
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "012203"
$newSheet.Range("C7").Value = "中加消费优选混合C"
$newSheet.Range("D7").Value = "0.95"
$newSheet.Range("E7").Value = "88.36"
$newSheet.Range("F7").Value = "3.91"
$newSheet.Range("G7").Value = "0.0371"
$newSheet.Range("H7").Value = 6
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "012522"
$newSheet.Range("C8").Value = "英大稳固增强核心一年持有混合C"
$newSheet.Range("D8").Value = "1.05"
$newSheet.Range("E8").Value = "23.17"
$newSheet.Range("F8").Value = "2.14"
$newSheet.Range("G8").Value = "0.0225"
$newSheet.Range("H8").Value = 1
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "004436"
$newSheet.Range("C9").Value = "汇添富年年泰定期开放混合A"
$newSheet.Range("D9").Value = "1.58"
$newSheet.Range("E9").Value = "24.12"
$newSheet.Range("F9").Value = "1.34"
$newSheet.Range("G9").Value = "0.0212"
$newSheet.Range("H9").Value = 7
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "008033"
$newSheet.Range("C10").Value = "中加科盈混合A"
$newSheet.Range("D10").Value = "1.25"
$newSheet.Range("E10").Value = "23.06"
$newSheet.Range("F10").Value = "1.46"
$newSheet.Range("G10").Value = "0.0182"
$newSheet.Range("H10").Value = 2
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "012521"
$newSheet.Range("C11").Value = "英大稳固增强核心一年持有混合A"
$newSheet.Range("D11").Value = "0.63"
$newSheet.Range("E11").Value = "23.17"
$newSheet.Range("F11").Value = "2.14"
$newSheet.Range("G11").Value = "0.0135"
$newSheet.Range("H11").Value = 1
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "004437"
$newSheet.Range("C12").Value = "汇添富年年泰定期开放混合C"
$newSheet.Range("D12").Value = "0.14"
$newSheet.Range("E12").Value = "24.12"
$newSheet.Range("F12").Value = "1.34"
$newSheet.Range("G12").Value = "0.0019"
$newSheet.Range("H12").Value = 7
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "001608"
$newSheet.Range("C13").Value = "英大策略优选混合C"
$newSheet.Range("D13").Value = "0.02"
$newSheet.Range("E13").Value = "93.12"
$newSheet.Range("F13").Value = "6.80"
$newSheet.Range("G13").Value = "0.0014"
$newSheet.Range("H13").Value = 4
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "008034"
$newSheet.Range("C14").Value = "中加科盈混合C"
$newSheet.Range("D14").Value = "0.07"
$newSheet.Range("E14").Value = "23.06"
$newSheet.Range("F14").Value = "1.46"
$newSheet.Range("G14").Value = "0.0010"
$newSheet.Range("H14").Value = 2

# --- 2. Insert the new 2022-Q4 row at the top of the "总计" data table ----
$summarySheet.Rows.Item(2).Insert()
$summarySheet.Range("B2:D2").ClearFormats()

# Give the new index cell (column A) the same style used by the other
# index cells in that column.
$summarySheet.Range("A3").Copy()
$summarySheet.Range("A2").PasteSpecial(-4122)

$summarySheet.Range("B2").Value = "2022-Q4"
$summarySheet.Range("C2").Value = 13
$summarySheet.Range("D2").Value = 1.11

# Column A holds a plain 0-based row counter; renumber the whole table
# (the insert above only shifted the rows down, it left the stale index
# values behind).
$summarySheet.Range("A2").Value = 0
$summarySheet.Range("A3").Value = 1
$summarySheet.Range("A4").Value = 2
$summarySheet.Range("A5").Value = 3
$summarySheet.Range("A6").Value = 4
$summarySheet.Range("A7").Value = 5
$summarySheet.Range("A8").Value = 6
$summarySheet.Range("A9").Value = 7
$summarySheet.Range("A10").Value = 8
